$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = [double]"3"
$ws.Range("F2").Value = [double]"1"
$ws.Range("G2").Value = [double]"2.276052666666667"
$ws.Range("H2").Value = [double]"6.828158"
$ws.Range("I2").Value = [double]"0.005247614157263819"
$ws.Range("J2").Value = [double]"0.005247614157263819"
$ws.Range("K2").Value = [double]"3"
$ws.Range("L2").Value = [double]"1"
$ws.Range("M2").Value = [double]"0.4702473333333333"
$ws.Range("N2").Value = [double]"1.410742"
$ws.Range("O2").Value = [double]"0.009034922268422819"
$ws.Range("P2").Value = [double]"0.009034922268422819"
$ws.Range("Q2").Value = [double]"1.070307697026222"
$ws.Range("R2").Value = [double]"9.632769273236001"
$ws.Range("S2").Value = [double]"4.741178600555372E-05"
$ws.Range("T2").Value = [double]"4.741178600555373E-05"

$ws.Range("E3").Value = [double]"3"
$ws.Range("F3").Value = [double]"1"
$ws.Range("G3").Value = [double]"2.276052666666667"
$ws.Range("H3").Value = [double]"6.828158"
$ws.Range("I3").Value = [double]"0.005247614157263819"
$ws.Range("J3").Value = [double]"0.005247614157263819"
$ws.Range("N3").Value = [double]"0.9584440000000001"
$ws.Range("O3").Value = [double]"0.006138235792679485"
$ws.Range("P3").Value = [double]"0.006138235792679485"
$ws.Range("Q3").Value = [double]"0.7271563406835556"
$ws.Range("R3").Value = [double]"6.544407066152001"
$ws.Range("S3").Value = [double]"3.221109304628837E-05"
$ws.Range("T3").Value = [double]"3.221109304628837E-05"

$ws.Range("E4").Value = [double]"3"
$ws.Range("F4").Value = [double]"1"
$ws.Range("G4").Value = [double]"2.276052666666667"
$ws.Range("H4").Value = [double]"6.828158"
$ws.Range("I4").Value = [double]"0.005247614157263819"
$ws.Range("J4").Value = [double]"0.005247614157263819"
$ws.Range("M4").Value = [double]"1.047307"
$ws.Range("N4").Value = [double]"3.141921"
$ws.Range("O4").Value = [double]"0.02012204358311108"
$ws.Range("P4").Value = [double]"0.02012204358311108"
$ws.Range("Q4").Value = [double]"2.383725890168667"
$ws.Range("R4").Value = [double]"21.453533011518"
$ws.Range("S4").Value = [double]"0.0001055927207798133"
$ws.Range("T4").Value = [double]"0.0001055927207798133"

$ws.Range("E5").Value = [double]"3"
$ws.Range("F5").Value = [double]"1"
$ws.Range("G5").Value = [double]"2.276052666666667"
$ws.Range("H5").Value = [double]"6.828158"
$ws.Range("I5").Value = [double]"0.005247614157263819"
$ws.Range("J5").Value = [double]"0.005247614157263819"
$ws.Range("M5").Value = [double]"50.21070966666667"
$ws.Range("N5").Value = [double]"150.632129"
$ws.Range("O5").Value = [double]"0.9647047983557866"
$ws.Range("P5").Value = [double]"0.9647047983557866"
$ws.Range("Q5").Value = [double]"114.2822196320425"
$ws.Range("R5").Value = [double]"1028.539976688382"
$ws.Range("S5").Value = [double]"0.005062398557432163"
$ws.Range("T5").Value = [double]"0.005062398557432164"

$ws.Range("I6").Value = [double]"0.1062533062835484"
$ws.Range("J6").Value = [double]"0.1062533062835484"
$ws.Range("K6").Value = [double]"3"
$ws.Range("L6").Value = [double]"1"
$ws.Range("M6").Value = [double]"0.4702473333333333"
$ws.Range("N6").Value = [double]"1.410742"
$ws.Range("O6").Value = [double]"0.009034922268422819"
$ws.Range("P6").Value = [double]"0.009034922268422819"
$ws.Range("Q6").Value = [double]"21.67151168924045"
$ws.Range("R6").Value = [double]"195.043605203164"
$ws.Range("S6").Value = [double]"0.0009599903630347816"
$ws.Range("T6").Value = [double]"0.0009599903630347817"

$ws.Range("I7").Value = [double]"0.1062533062835484"
$ws.Range("J7").Value = [double]"0.1062533062835484"
$ws.Range("N7").Value = [double]"0.9584440000000001"
$ws.Range("O7").Value = [double]"0.006138235792679485"
$ws.Range("P7").Value = [double]"0.006138235792679485"
$ws.Range("S7").Value = [double]"0.0006522078477202127"
$ws.Range("T7").Value = [double]"0.0006522078477202128"

$ws.Range("I8").Value = [double]"0.1062533062835484"
$ws.Range("J8").Value = [double]"0.1062533062835484"
$ws.Range("M8").Value = [double]"1.047307"
$ws.Range("N8").Value = [double]"3.141921"
$ws.Range("O8").Value = [double]"0.02012204358311108"
$ws.Range("P8").Value = [double]"0.02012204358311108"
$ws.Range("Q8").Value = [double]"48.26550685963134"
$ws.Range("R8").Value = [double]"434.389561736682"
$ws.Range("S8").Value = [double]"0.002138033659887211"
$ws.Range("T8").Value = [double]"0.002138033659887211"

$ws.Range("I9").Value = [double]"0.1062533062835484"
$ws.Range("J9").Value = [double]"0.1062533062835484"
$ws.Range("M9").Value = [double]"50.21070966666667"
$ws.Range("N9").Value = [double]"150.632129"
$ws.Range("O9").Value = [double]"0.9647047983557866"
$ws.Range("P9").Value = [double]"0.9647047983557866"
$ws.Range("Q9").Value = [double]"2313.977994841492"
$ws.Range("R9").Value = [double]"20825.80195357342"
$ws.Range("S9").Value = [double]"0.1025030744129062"
$ws.Range("T9").Value = [double]"0.1025030744129062"

$ws.Range("G10").Value = [double]"41.187613"
$ws.Range("H10").Value = [double]"123.562839"
$ws.Range("I10").Value = [double]"0.09496120377532416"
$ws.Range("J10").Value = [double]"0.09496120377532417"
$ws.Range("K10").Value = [double]"3"
$ws.Range("L10").Value = [double]"1"
$ws.Range("M10").Value = [double]"0.4702473333333333"
$ws.Range("N10").Value = [double]"1.410742"
$ws.Range("O10").Value = [double]"0.009034922268422819"
$ws.Range("P10").Value = [double]"0.009034922268422819"
$ws.Range("Q10").Value = [double]"19.36836517961533"
$ws.Range("R10").Value = [double]"174.315286616538"
$ws.Range("S10").Value = [double]"0.0008579670946259134"
$ws.Range("T10").Value = [double]"0.0008579670946259135"

$ws.Range("G11").Value = [double]"41.187613"
$ws.Range("H11").Value = [double]"123.562839"
$ws.Range("I11").Value = [double]"0.09496120377532416"
$ws.Range("J11").Value = [double]"0.09496120377532417"
$ws.Range("N11").Value = [double]"0.9584440000000001"
$ws.Range("O11").Value = [double]"0.006138235792679485"
$ws.Range("P11").Value = [double]"0.006138235792679485"
$ws.Range("Q11").Value = [double]"13.15867351805733"
$ws.Range("R11").Value = [double]"118.428061662516"
$ws.Range("S11").Value = [double]"0.0005828942599296251"
$ws.Range("T11").Value = [double]"0.0005828942599296251"

$ws.Range("G12").Value = [double]"41.187613"
$ws.Range("H12").Value = [double]"123.562839"
$ws.Range("I12").Value = [double]"0.09496120377532416"
$ws.Range("J12").Value = [double]"0.09496120377532417"
$ws.Range("M12").Value = [double]"1.047307"
$ws.Range("N12").Value = [double]"3.141921"
$ws.Range("O12").Value = [double]"0.02012204358311108"
$ws.Range("P12").Value = [double]"0.02012204358311108"
$ws.Range("Q12").Value = [double]"43.136075408191"
$ws.Range("R12").Value = [double]"388.224678673719"
$ws.Range("S12").Value = [double]"0.001910813481071765"
$ws.Range("T12").Value = [double]"0.001910813481071766"

$ws.Range("G13").Value = [double]"41.187613"
$ws.Range("H13").Value = [double]"123.562839"
$ws.Range("I13").Value = [double]"0.09496120377532416"
$ws.Range("J13").Value = [double]"0.09496120377532417"
$ws.Range("M13").Value = [double]"50.21070966666667"
$ws.Range("N13").Value = [double]"150.632129"
$ws.Range("O13").Value = [double]"0.9647047983557866"
$ws.Range("P13").Value = [double]"0.9647047983557866"
$ws.Range("Q13").Value = [double]"2068.059278206026"
$ws.Range("R13").Value = [double]"18612.53350385423"
$ws.Range("S13").Value = [double]"0.09160952893969686"
$ws.Range("T13").Value = [double]"0.09160952893969687"

$ws.Range("G14").Value = [double]"344.1819356666667"
$ws.Range("H14").Value = [double]"1032.545807"
$ws.Range("I14").Value = [double]"0.7935378757838636"
$ws.Range("J14").Value = [double]"0.7935378757838637"
$ws.Range("K14").Value = [double]"3"
$ws.Range("L14").Value = [double]"1"
$ws.Range("M14").Value = [double]"0.4702473333333333"
$ws.Range("N14").Value = [double]"1.410742"
$ws.Range("O14").Value = [double]"0.009034922268422819"
$ws.Range("P14").Value = [double]"0.009034922268422819"
$ws.Range("Q14").Value = [double]"161.8506374287549"
$ws.Range("R14").Value = [double]"1456.655736858794"
$ws.Range("S14").Value = [double]"0.00716955302475657"
$ws.Range("T14").Value = [double]"0.007169553024756571"

$ws.Range("G15").Value = [double]"344.1819356666667"
$ws.Range("H15").Value = [double]"1032.545807"
$ws.Range("I15").Value = [double]"0.7935378757838636"
$ws.Range("J15").Value = [double]"0.7935378757838637"
$ws.Range("N15").Value = [double]"0.9584440000000001"
$ws.Range("O15").Value = [double]"0.006138235792679485"
$ws.Range("P15").Value = [double]"0.006138235792679485"
$ws.Range("Q15").Value = [double]"109.9597037160342"
$ws.Range("R15").Value = [double]"989.6373334443081"
$ws.Range("S15").Value = [double]"0.004870922591983359"
$ws.Range("T15").Value = [double]"0.00487092259198336"

$ws.Range("G16").Value = [double]"344.1819356666667"
$ws.Range("H16").Value = [double]"1032.545807"
$ws.Range("I16").Value = [double]"0.7935378757838636"
$ws.Range("J16").Value = [double]"0.7935378757838637"
$ws.Range("M16").Value = [double]"1.047307"
$ws.Range("N16").Value = [double]"3.141921"
$ws.Range("O16").Value = [double]"0.02012204358311108"
$ws.Range("P16").Value = [double]"0.02012204358311108"
$ws.Range("Q16").Value = [double]"360.4641504972496"
$ws.Range("R16").Value = [double]"3244.177354475247"
$ws.Range("S16").Value = [double]"0.01596760372137229"
$ws.Range("T16").Value = [double]"0.0159676037213723"

$ws.Range("G17").Value = [double]"344.1819356666667"
$ws.Range("H17").Value = [double]"1032.545807"
$ws.Range("I17").Value = [double]"0.7935378757838636"
$ws.Range("J17").Value = [double]"0.7935378757838637"
$ws.Range("M17").Value = [double]"50.21070966666667"
$ws.Range("N17").Value = [double]"150.632129"
$ws.Range("O17").Value = [double]"0.9647047983557866"
$ws.Range("P17").Value = [double]"0.9647047983557866"
$ws.Range("Q17").Value = [double]"17281.61924427035"
$ws.Range("R17").Value = [double]"155534.5731984331"
$ws.Range("S17").Value = [double]"0.7655297964457514"
$ws.Range("T17").Value = [double]"0.7655297964457515"
